$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Swap the values of D6 and E6: D6 becomes "Admin", E6 becomes "admin123"
$ws.Range("D6").Value = "Admin"
$ws.Range("E6").Value = "admin123"

# Update the selection on the DATA sheet to E6
$ws.Activate()
$ws.Range("E6").Select()
